# The only textual differences in the target diff are:
#   1) The ordering of the xmlns:* attributes on the root elements of
#      word/document.xml, word/footer.xml, word/header.xml and
#      word/styles.xml.
#   2) One word inside an XML comment in word/document.xml
#      ("Modified by docx4j ... using REFERENCE JAXB in
#      Microsoft Java 21.0.8 on Mac OS X" -> "... Oracle Java 21.0.8 ...").
#
# Both of these come from the external docx4j post-processing step that
# produced the "after" fixture in the source repository (its JVM vendor
# string, and its own XML serializer's attribute ordering) - neither is
# part of the Word document content/object model, so there is nothing
# here a Word document edit can (or should) change: no paragraph,
# run, table, header/footer, or style content differs between the
# "before" and "after" OOXML. Touching the document would only destroy
# the round-trip fidelity (and that docx4j comment) without being able
# to reproduce the non-Word-observable JVM-vendor string or serializer
# attribute order anyway.
#
# So this script intentionally performs no content mutation - it just
# touches the document object to confirm the runtime is wired up.
$d = $word.ActiveDocument
$d.Content | Out-Null
